$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.497.41'
$ws.Range("E2").Value = '  -2.76%  '
$ws.Range("D3").Value = '2.247.04'
$ws.Range("E3").Value = '  -3.79%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.11'
$ws.Range("E5").Value = '  -1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.634'
$ws.Range("E6").Value = '  -4.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.91'
$ws.Range("E7").Value = '  -3.05%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  -5.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0999'
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.09'
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '36.84'
$ws.Range("E12").Value = '  +14.31%  '
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.76'
$ws.Range("E14").Value = '  -5.19%  '
$ws.Range("D15").Value = '2.580.85'
$ws.Range("E15").Value = '  -3.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.16'
$ws.Range("E16").Value = '  -5.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.859'
$ws.Range("E17").Value = '  -4.03%  '
$ws.Range("D18").Value = '2.247.72'
$ws.Range("E18").Value = '  -3.74%  '
$ws.Range("D19").Value = '42.271.82'
$ws.Range("E19").Value = '  -3.13%  '
$ws.Range("D20").Value = '0.0₃0974'
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.28'
$ws.Range("E21").Value = '  -4.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.41'
$ws.Range("E22").Value = '  -5.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.79'
$ws.Range("E23").Value = '  -5.61%  '
$ws.Range("E24").Value = '  +3.51%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("E27").Value = '  -3.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.01'
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.49'
$ws.Range("E30").Value = '  -3.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.60'
$ws.Range("E31").Value = '  -6.74%  '
$ws.Range("E32").Value = '  -3.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.127'
$ws.Range("E33").Value = '  -5.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0724'
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.33'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.73'
$ws.Range("E36").Value = '  -6.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.69'
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.39'
$ws.Range("E38").Value = '  +20.04%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.29'
$ws.Range("E39").Value = '  -3.02%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0277'
$ws.Range("E40").Value = '  +2.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.95'
$ws.Range("E41").Value = '  -6.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.25'
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.28'
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.99'
$ws.Range("E44").Value = '  -12.68%  '
$ws.Range("E45").Value = '  -2.35%  '
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("B47").Value = 'SynthetixNetwork'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.64'
$ws.Range("E47").Value = '  +13.94%  '
$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.22'
$ws.Range("E49").Value = '  +10.04%  '
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("E51").Value = '  -2.63%  '
